{"js": "// Update the date and each two-digit \u00d7 two-digit multiplication\n// problem/answer cell to the new values.\n// Old -> New text pairs, in document order (date paragraph first,\n// then each table cell, left-to-right / top-to-bottom).\nconst pairs = [\n  [\"2024-02-11 Sunday\", \"2024-02-12 Monday\"],\n  [\"22\u00d786=1892\", \"22\u00d717=374\"],\n  [\"97\u00d791=8827\", \"44\u00d784=3696\"],\n  [\"93\u00d755=5115\", \"77\u00d781=6237\"],\n  [\"44\u00d755=2420\", \"20\u00d733=660\"],\n  [\"84\u00d738=3192\", \"56\u00d726=1456\"],\n  [\"40\u00d764=2560\", \"69\u00d716=1104\"],\n  [\"85\u00d788=7480\", \"44\u00d762=2728\"],\n  [\"11\u00d736=396\", \"86\u00d788=7568\"],\n  [\"13\u00d790=1170\", \"48\u00d765=3120\"],\n  [\"52\u00d738=1976\", \"26\u00d763=1638\"],\n  [\"50\u00d790=4500\", \"44\u00d764=2816\"],\n  [\"90\u00d751=4590\", \"24\u00d739=936\"],\n  [\"98\u00d767=6566\", \"19\u00d793=1767\"],\n  [\"87\u00d769=6003\", \"40\u00d777=3080\"],\n  [\"15\u00d777=1155\", \"17\u00d761=1037\"],\n  [\"24\u00d768=1632\", \"70\u00d774=5180\"],\n  [\"92\u00d724=2208\", \"17\u00d794=1598\"],\n  [\"94\u00d716=1504\", \"63\u00d729=1827\"],\n  [\"24\u00d795=2280\", \"61\u00d717=1037\"],\n  [\"71\u00d724=1704\", \"96\u00d795=9120\"],\n  [\"22\u00d756=1232\", \"52\u00d746=2392\"],\n  [\"89\u00d746=4094\", \"75\u00d755=4125\"],\n  [\"99\u00d729=2871\", \"88\u00d731=2728\"],\n  [\"54\u00d727=1458\", \"84\u00d766=5544\"],\n  [\"94\u00d743=4042\", \"50\u00d795=4750\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and each two-digit x two-digit multiplication\n# problem/answer table cell to the new values using Find & Replace.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-02-11 Sunday\", \"2024-02-12 Monday\"),\n    @(\"22\u00d786=1892\", \"22\u00d717=374\"),\n    @(\"97\u00d791=8827\", \"44\u00d784=3696\"),\n    @(\"93\u00d755=5115\", \"77\u00d781=6237\"),\n    @(\"44\u00d755=2420\", \"20\u00d733=660\"),\n    @(\"84\u00d738=3192\", \"56\u00d726=1456\"),\n    @(\"40\u00d764=2560\", \"69\u00d716=1104\"),\n    @(\"85\u00d788=7480\", \"44\u00d762=2728\"),\n    @(\"11\u00d736=396\", \"86\u00d788=7568\"),\n    @(\"13\u00d790=1170\", \"48\u00d765=3120\"),\n    @(\"52\u00d738=1976\", \"26\u00d763=1638\"),\n    @(\"50\u00d790=4500\", \"44\u00d764=2816\"),\n    @(\"90\u00d751=4590\", \"24\u00d739=936\"),\n    @(\"98\u00d767=6566\", \"19\u00d793=1767\"),\n    @(\"87\u00d769=6003\", \"40\u00d777=3080\"),\n    @(\"15\u00d777=1155\", \"17\u00d761=1037\"),\n    @(\"24\u00d768=1632\", \"70\u00d774=5180\"),\n    @(\"92\u00d724=2208\", \"17\u00d794=1598\"),\n    @(\"94\u00d716=1504\", \"63\u00d729=1827\"),\n    @(\"24\u00d795=2280\", \"61\u00d717=1037\"),\n    @(\"71\u00d724=1704\", \"96\u00d795=9120\"),\n    @(\"22\u00d756=1232\", \"52\u00d746=2392\"),\n    @(\"89\u00d746=4094\", \"75\u00d755=4125\"),\n    @(\"99\u00d729=2871\", \"88\u00d731=2728\"),\n    @(\"54\u00d727=1458\", \"84\u00d766=5544\"),\n    @(\"94\u00d743=4042\", \"50\u00d795=4750\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
